# Update of league bases - 14-05-2024 20:19
# Costa Rica Primera Division workbook update:
#  - Several pairs/triples of existing rows had their match data corrected
#    (the "id" in column A is left untouched, everything else is swapped
#    between the affected rows).
#  - Six brand new match rows are appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows([int]$r1, [int]$r2) {
    $range1 = "B$($r1):AB$($r1)"
    $range2 = "B$($r2):AB$($r2)"
    $v1 = $ws.Range($range1).Value2
    $v2 = $ws.Range($range2).Value2
    $ws.Range($range1).Value = $v2
    $ws.Range($range2).Value = $v1
}

# Rotate the content (columns B..AB) of the given rows: the content that was
# in rows[1] ends up in rows[0], rows[2] -> rows[1], ..., rows[0] -> rows[n-1].
# (i.e. each row receives the content that used to belong to the next row in
# the list, wrapping around at the end)
function Rotate-Rows([int[]]$rows) {
    $vals = @()
    foreach ($r in $rows) {
        $vals += , ($ws.Range("B$($r):AB$($r)").Value2)
    }
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $src = $vals[($i + 1) % $rows.Length]
        $dstRow = $rows[$i]
        $ws.Range("B$($dstRow):AB$($dstRow)").Value = $src
    }
}

# --- Corrections to existing rows (id column A is unaffected) -------------
Swap-Rows 38 39
Swap-Rows 91 92
Rotate-Rows @(129, 130, 131)
Swap-Rows 200 201
Swap-Rows 249 250

$excel.CutCopyMode = $false

# --- Append six brand-new match rows (266-271) -----------------------------
$newRows = @(
    @{ Row=266; A=264; B=8162894; C="Costa Rica Primera Division"; D=45423.79166666666; E="Puntarenas";          F="Cartagines";         G=1; H=1; I="D"; J=2.3;   K=3.3;  L=2.75;  M=2.6;  N=3.2;   O=2.5;   P=0;     Q=1.975; R=1.825; S=2.5;  T=1.85;  U=1.95;  V=-1; W=2.2; X=-1; Y=0;      Z=0;      AA=-1;     AB=0.95 },
    @{ Row=267; A=265; B=8162891; C="Costa Rica Primera Division"; D=45424.75;          E="Deportivo Saprissa";  F="Santos de Gupiles";  G=3; H=1; I="H"; J=1.166; K=6.5;  L=13;    M=1.125; N=8;     O=15;    P=-2.25; Q=1.825; R=1.975; S=3.5;  T=1.975; U=1.825; V=0.125; W=-1; X=-1; Y=-0.5;   Z=0.4875; AA=0.9750000000000001; AB=-1 },
    @{ Row=268; A=266; B=8162893; C="Costa Rica Primera Division"; D=45424.75;          E="AD Grecia";           F="AD San Carlos";      G=2; H=2; I="D"; J=5;     K=4;    L=1.533; M=4.2;  N=4.2;   O=1.6;   P=1;     Q=1.775; R=2.025; S=3;    T=1.925; U=1.875; V=-1; W=3.2; X=-1; Y=0.7749999999999999; Z=-1; AA=0.925;  AB=-1 },
    @{ Row=269; A=267; B=8203655; C="Costa Rica Primera Division"; D=45424.75;          E="Municipal Perez Zeledon"; F="Municipal Liberia"; G=0; H=3; I="A"; J=3.3; K=3.5; L=2;     M=2.9;  N=3.4;   O=2.2;   P=0.25;  Q=1.825; R=1.975; S=2.75; T=2;     U=1.8;   V=-1; W=-1; X=1.2; Y=-1;   Z=0.9750000000000001; AA=0.5; AB=-0.5 },
    @{ Row=270; A=268; B=8162895; C="Costa Rica Primera Division"; D=45424.75;          E="Sporting San Jose";   F="Herediano";          G=1; H=1; I="D"; J=3.6;   K=3.5;  L=1.833; M=4.5;  N=3.8;   O=1.571; P=0.75;  Q=2.025; R=1.775; S=2.75; T=1.975; U=1.825; V=-1; W=2.8; X=-1; Y=1.025;  Z=-1;     AA=-1;     AB=0.825 },
    @{ Row=271; A=269; B=8162892; C="Costa Rica Primera Division"; D=45424.75;          E="Alajuelense";         F="AD Guanacasteca";    G=5; H=0; I="H"; J=1.25;  K=5;    L=10;    M=1.3;  N=4.75;  O=8;     P=-1.5;  Q=1.9;   R=1.9;   S=3;    T=1.9;   U=1.9;   V=0.3; W=-1; X=-1; Y=0.8999999999999999; Z=-1; AA=0.8999999999999999; AB=-1 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $prev = $r - 1
    # Copy the formatting (styles) of the row above, then fill in the values.
    $ws.Range("A$($prev):AB$($prev)").Copy()
    $ws.Range("A$($r):AB$($r)").PasteSpecial(-4122) # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Range("A$r").Value = $nr.A
    $ws.Range("B$r").Value = $nr.B
    $ws.Range("C$r").Value = $nr.C
    $ws.Range("D$r").Value = $nr.D
    $ws.Range("E$r").Value = $nr.E
    $ws.Range("F$r").Value = $nr.F
    $ws.Range("G$r").Value = $nr.G
    $ws.Range("H$r").Value = $nr.H
    $ws.Range("I$r").Value = $nr.I
    $ws.Range("J$r").Value = $nr.J
    $ws.Range("K$r").Value = $nr.K
    $ws.Range("L$r").Value = $nr.L
    $ws.Range("M$r").Value = $nr.M
    $ws.Range("N$r").Value = $nr.N
    $ws.Range("O$r").Value = $nr.O
    $ws.Range("P$r").Value = $nr.P
    $ws.Range("Q$r").Value = $nr.Q
    $ws.Range("R$r").Value = $nr.R
    $ws.Range("S$r").Value = $nr.S
    $ws.Range("T$r").Value = $nr.T
    $ws.Range("U$r").Value = $nr.U
    $ws.Range("V$r").Value = $nr.V
    $ws.Range("W$r").Value = $nr.W
    $ws.Range("X$r").Value = $nr.X
    $ws.Range("Y$r").Value = $nr.Y
    $ws.Range("Z$r").Value = $nr.Z
    $ws.Range("AA$r").Value = $nr.AA
    $ws.Range("AB$r").Value = $nr.AB
}

Write-Host "Applied league base update."
